$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 23 (pushing existing rows 23-71 down to 25-73).
# Using EntireRow.Insert() on a 2-row range inserts 2 rows at once, shifting
# cell formatting down with them (so D23/D24 inherit the date number format
# from the row above, same as the surrounding rows).
$ws.Range("A23:A24").EntireRow.Insert()

# New row 23: Membrillo Champion "Especial"
$ws.Range("A23").Value = 9
$ws.Range("B23").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C23").Value = "Metropolitana"
$ws.Range("D23").Value = 45054
$ws.Range("E23").Value = 13
$ws.Range("F23").Value = "Fruta"
$ws.Range("G23").Value = 100104
$ws.Range("H23").Value = "Frutos de pepita"
$ws.Range("I23").Value = 100104003
$ws.Range("J23").Value = "Membrillo"
$ws.Range("K23").Value = "Champion"
$ws.Range("L23").Value = "Especial"
$ws.Range("M23").Value = 250
$ws.Range("N23").Value = 10500
$ws.Range("O23").Value = 10500
$ws.Range("P23").Value = 10500
$ws.Range("Q23").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R23").Value = "Región de O'Higgins"
$ws.Range("S23").Value = 700
$ws.Range("T23").Value = 15

# New row 24: Membrillo Champion "Primera"
$ws.Range("A24").Value = 9
$ws.Range("B24").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C24").Value = "Metropolitana"
$ws.Range("D24").Value = 45054
$ws.Range("E24").Value = 13
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100104
$ws.Range("H24").Value = "Frutos de pepita"
$ws.Range("I24").Value = 100104003
$ws.Range("J24").Value = "Membrillo"
$ws.Range("K24").Value = "Champion"
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 220
$ws.Range("N24").Value = 7500
$ws.Range("O24").Value = 7500
$ws.Range("P24").Value = 7500
$ws.Range("Q24").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R24").Value = "Región de O'Higgins"
$ws.Range("S24").Value = 500
$ws.Range("T24").Value = 15
